$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit adds one new weekly price record for "Poroto verde" at
# Mercado Mayorista Lo Valledor de Santiago. It is inserted as row 1250,
# pushing the existing rows 1250-1295 down to 1251-1296 (dimension grows
# from A1:R1295 to A1:R1296).
$ws.Rows(1250).Insert()

$ws.Cells.Item(1250, 1).Value = 6
$ws.Cells.Item(1250, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1250, 3).Value = "Metropolitana"
$ws.Cells.Item(1250, 4).Value = 45075
$ws.Cells.Item(1250, 5).Value = 13
$ws.Cells.Item(1250, 6).Value = 100112031
$ws.Cells.Item(1250, 7).Value = "Poroto verde"
$ws.Cells.Item(1250, 8).Value = "Magnum"
$ws.Cells.Item(1250, 9).Value = "Primera"
$ws.Cells.Item(1250, 10).Value = 320
$ws.Cells.Item(1250, 11).Value = 16000
$ws.Cells.Item(1250, 12).Value = 18000
$ws.Cells.Item(1250, 13).Value = 16625
$ws.Cells.Item(1250, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(1250, 15).Value = "Perú"
$ws.Cells.Item(1250, 16).Value = 665
$ws.Cells.Item(1250, 17).Value = 25
$ws.Cells.Item(1250, 18).Value = "Hortaliza"
